$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 85) with the latest Argent price entry, matching the
# plain-text format used by the existing data rows (dates/values are stored
# as literal text, not auto-converted to date serials / numbers). A leading
# apostrophe forces text entry; resetting the style afterwards keeps the
# cell on the default (unstyled) format, same as the preceding data rows.
$cellA = $ws.Cells.Item(85, 1)
$cellA.Value = "'2025-01-22"
$cellA.Style = "Normal"

$cellB = $ws.Cells.Item(85, 2)
$cellB.Value = "'42.6"
$cellB.Style = "Normal"
